# Weekly market-price refresh for the Leve profit tables (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H,I,J,K,L,M,N) for the
# rows whose source prices changed; these are static snapshot values (no formulas in the sheet).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 125: Body over Mind
$ws.Range("H125").Value = 2625.6428
$ws.Range("J125").Value = 3932.6667
$ws.Range("L125").Value = 35394.0003
$ws.Range("N125").Value = -40314.0003

# Row 127: Liquid Competence
$ws.Range("H127").Value = 461.16666
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

# Row 129: Practical Command
$ws.Range("H129").Value = 2999.6667
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 2999.6667
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 8999.000100000001
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -18999.0001

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 3142.3462
$ws.Range("I137").Value = 2546.077
$ws.Range("K137").Value = 7638.231000000001
$ws.Range("M137").Value = -5088.231000000001

# Row 141: Remedy for Reason
$ws.Range("H141").Value = 7000
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 7000
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 21000
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -31360

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 5952.8945
$ws.Range("I32").Value = 4604.2256
$ws.Range("K32").Value = 4604.2256
$ws.Range("M32").Value = -4317.2256

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 3944.6956
$ws.Range("I132").Value = 2700.6667
$ws.Range("J132").Value = 17007
$ws.Range("K132").Value = 8102.000100000001
$ws.Range("L132").Value = 51021
$ws.Range("M132").Value = -5572.000100000001
$ws.Range("N132").Value = -56081

$ws = $wb.Worksheets.Item("BSM")
# Row 21: Awl or Nothing
$ws.Range("H21").Value = 54982.25
$ws.Range("J21").Value = 54982.25
$ws.Range("L21").Value = 54982.25
$ws.Range("N21").Value = -55454.25

# Row 105: Ingot to Wing It
$ws.Range("H105").Value = 25461
$ws.Range("I105").Value = 42348.6
$ws.Range("J105").Value = 13398.429
$ws.Range("K105").Value = 42348.6
$ws.Range("L105").Value = 13398.429
$ws.Range("M105").Value = -40601.6
$ws.Range("N105").Value = -16892.429

$ws = $wb.Worksheets.Item("CRP")
# Row 21: Nightmare on My Street
$ws.Range("H21").Value = 172.4
$ws.Range("J21").Value = 212.25
$ws.Range("L21").Value = 212.25
$ws.Range("N21").Value = -682.25

# Row 22: Driving Up the Wall
$ws.Range("H22").Value = 1654.8889
$ws.Range("I22").Value = 344.14285
$ws.Range("K22").Value = 344.14285
$ws.Range("M22").Value = 5.85714999999999

# Row 31: Wall Not Found
$ws.Range("H31").Value = 37611.09
$ws.Range("I31").Value = 5851.706
$ws.Range("J31").Value = 69370.47
$ws.Range("K31").Value = 5851.706
$ws.Range("L31").Value = 69370.47
$ws.Range("M31").Value = -5556.706
$ws.Range("N31").Value = -69960.47

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 37611.09
$ws.Range("I34").Value = 5851.706
$ws.Range("J34").Value = 69370.47
$ws.Range("K34").Value = 5851.706
$ws.Range("L34").Value = 69370.47
$ws.Range("M34").Value = -5649.706
$ws.Range("N34").Value = -69774.47

# Row 43: The Long Lance of the Law
$ws.Range("H43").Value = 30000
$ws.Range("J43").Value = 30000
$ws.Range("L43").Value = 30000
$ws.Range("N43").Value = -30368

# Row 96: Composition
$ws.Range("H96").Value = 19420.25
$ws.Range("J96").Value = 19420.25
$ws.Range("L96").Value = 19420.25
$ws.Range("N96").Value = -24912.25

# Row 99: O Pine
$ws.Range("H99").Value = 3622.2856
$ws.Range("I99").Value = 3438.7778
$ws.Range("K99").Value = 3438.7778
$ws.Range("M99").Value = -1940.7778

# Row 101: Everybody's Heard about the 'Berd
$ws.Range("H101").Value = 30000
$ws.Range("J101").Value = 30000
$ws.Range("L101").Value = 30000
$ws.Range("N101").Value = -36490

# Row 107: Built to Last
$ws.Range("H107").Value = 868.9487
$ws.Range("J107").Value = 917.6
$ws.Range("L107").Value = 917.6
$ws.Range("N107").Value = -4757.6

# Row 126: A Better Conductor
$ws.Range("H126").Value = 3622.2856
$ws.Range("I126").Value = 3438.7778
$ws.Range("K126").Value = 10316.3334
$ws.Range("M126").Value = -7846.3334

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 2722.742
$ws.Range("I132").Value = 1917.96
$ws.Range("J132").Value = 6076
$ws.Range("K132").Value = 5753.88
$ws.Range("L132").Value = 18228
$ws.Range("M132").Value = -3223.88
$ws.Range("N132").Value = -23288

$ws = $wb.Worksheets.Item("CUL")
# Row 113: Can't Eat Just One
$ws.Range("H113").Value = 1660.0555
$ws.Range("J113").Value = 1700.9231
$ws.Range("L113").Value = 5102.7693
$ws.Range("N113").Value = -9442.7693

# Row 131: The Mountain Steeped
$ws.Range("H131").Value = 7261839
$ws.Range("I131").Value = 17858028
$ws.Range("J131").Value = 5257154.5
$ws.Range("K131").Value = 53574084
$ws.Range("L131").Value = 15771463.5
$ws.Range("M131").Value = -53569044
$ws.Range("N131").Value = -15781543.5

$ws = $wb.Worksheets.Item("GSM")
# Row 44: Actually, It's Loyalty
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

# Row 62: The Goggles, They Do Naught
$ws.Range("H62").Value = 40768.8
$ws.Range("I62").Value = 39614.668
$ws.Range("K62").Value = 39614.668
$ws.Range("M62").Value = -38928.668

# Row 65: Peril Never Wore Safety Goggles (L)
$ws.Range("H65").Value = 40768.8
$ws.Range("I65").Value = 39614.668
$ws.Range("K65").Value = 118844.004
$ws.Range("M65").Value = -115412.004

# Row 102: Put the Metal to the Peddle
$ws.Range("H102").Value = 2658
$ws.Range("I102").Value = 1995.6842
$ws.Range("K102").Value = 1995.6842
$ws.Range("M102").Value = -373.6841999999999

# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 5647.974
$ws.Range("I122").Value = 4927.44
$ws.Range("J122").Value = 6934.643
$ws.Range("K122").Value = 14782.32
$ws.Range("L122").Value = 20803.929
$ws.Range("M122").Value = -12332.32
$ws.Range("N122").Value = -25703.929

# Row 126: Gold Rush Order
$ws.Range("H126").Value = 6041.3
$ws.Range("I126").Value = 3133
$ws.Range("K126").Value = 9399
$ws.Range("M126").Value = -6929

# Row 132: On Board for Lar
$ws.Range("H132").Value = 7169.9785
$ws.Range("I132").Value = 6617.816
$ws.Range("J132").Value = 9501.333000000001
$ws.Range("K132").Value = 19853.448
$ws.Range("L132").Value = 28503.999
$ws.Range("M132").Value = -17323.448
$ws.Range("N132").Value = -33563.999

# Row 135: Fan of the Foreign
$ws.Range("H135").Value = 69612.86
$ws.Range("J135").Value = 69612.86
$ws.Range("L135").Value = 69612.86
$ws.Range("N135").Value = -79752.86

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban
$ws.Range("H7").Value = 8345.132
$ws.Range("I7").Value = 6721.5415
$ws.Range("K7").Value = 6721.5415
$ws.Range("M7").Value = -6609.5415

# Row 40: Best Served Toad
$ws.Range("H40").Value = 7068.383
$ws.Range("I40").Value = 6836.528
$ws.Range("K40").Value = 6836.528
$ws.Range("M40").Value = -6700.528

# Row 46: Supply Side Logic
$ws.Range("H46").Value = 3496.3704
$ws.Range("I46").Value = 2081.9092
$ws.Range("J46").Value = 4468.8125
$ws.Range("K46").Value = 2081.9092
$ws.Range("L46").Value = 4468.8125
$ws.Range("M46").Value = -1893.9092
$ws.Range("N46").Value = -4844.8125

# Row 68: You Could Say It's a Moving Target
$ws.Range("H68").Value = 4789.25
$ws.Range("I68").Value = 4069.1428
$ws.Range("K68").Value = 4069.1428
$ws.Range("M68").Value = -3320.1428

# Row 71: They Call It Bloody Mary (L)
$ws.Range("H71").Value = 4789.25
$ws.Range("I71").Value = 4069.1428
$ws.Range("K71").Value = 20345.714
$ws.Range("M71").Value = -16601.714

# Row 122: Hell on Leather
$ws.Range("H122").Value = 114441.914
$ws.Range("I122").Value = 144856.67
$ws.Range("K122").Value = 434570.01
$ws.Range("M122").Value = -432120.01

# Row 126: Battered Books
$ws.Range("H126").Value = 8345.132
$ws.Range("I126").Value = 6721.5415
$ws.Range("K126").Value = 20164.6245
$ws.Range("M126").Value = -17694.6245

# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 7459.6045
$ws.Range("I132").Value = 6809.8
$ws.Range("J132").Value = 8362.111000000001
$ws.Range("K132").Value = 20429.4
$ws.Range("L132").Value = 25086.333
$ws.Range("M132").Value = -17899.4
$ws.Range("N132").Value = -30146.333

# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 6856.9707
$ws.Range("I136").Value = 2697.44
$ws.Range("K136").Value = 8092.32
$ws.Range("M136").Value = -5542.32

$ws = $wb.Worksheets.Item("WVR")
# Row 38: By the Seat of the Pants
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()

# Row 107: Flax Wax
$ws.Range("H107").Value = 1212.8572
$ws.Range("I107").Value = 998.3333
$ws.Range("K107").Value = 2994.9999
$ws.Range("M107").Value = -1074.9999

# Row 126: A Polished Purchase
$ws.Range("H126").Value = 1749.0435
$ws.Range("I126").Value = 1361.3125
$ws.Range("K126").Value = 4083.9375
$ws.Range("M126").Value = -1613.9375

# Row 132: Comfy Cabins
$ws.Range("H132").Value = 4822.727
$ws.Range("I132").Value = 3005.647
$ws.Range("K132").Value = 9016.940999999999
$ws.Range("M132").Value = -6486.940999999999
